$d = $word.ActiveDocument

# 1. Insert a new Title paragraph at the very start of the document body:
#    "Thermoregulation notebook"
$range = $d.Range(0, 0)
$range.InsertBefore("Thermoregulation notebook`r")
$d.Paragraphs(1).Style = "Title"

# 2. "You may generate the Excel sheet" -> "You may work the Excel sheet"
#    and "answers to the question below" -> "answers to the questions below"
$d.Content.Find.Execute(
    "You may generate the Excel sheet as part of your group, but your analysis and answers to the question below should be your own work and in your own words.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You may work the Excel sheet as part of your group, but your analysis and answers to the questions below should be your own work and in your own words.",
    2) | Out-Null

# 3. "In your Excel file, put in formulae to calculate ... Calculate the mean half-life per treatment."
#    -> "In your Excel file, I have entered the formulae to calculate ... Use Excel formulae to calculate the mean half-life per treatment."
$d.Content.Find.Execute(
    "In your Excel file, put in formulae to calculate the half-life for each 30 second segment of all of the treatments (it’s one formula, then copy and paste). Calculate the mean half-life per treatment. Put those values below for the treatments you were able to complete:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In your Excel file, I have entered the formulae to calculate the half-life for each 30 second segment of all of the treatments (it’s one formula, then copy and paste). Use Excel formulae to calculate the mean half-life per treatment. Put those values below for the treatments you were able to complete:",
    2) | Out-Null

# 4. "treatment. Explain in your own words what happened to the mouse when you fan it using the"
#    -> "treatment. Explain in your own words what happened to the mouse when you fan it; use the"
#    Keep the edit confined to the single run that originally holds this text, so
#    surrounding runs (the " ", the curly-quoted "wind chill" run, etc.) are left alone.
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Text = "happened to the mouse when you fan it using the"
$r1.Find.Forward = $true
if ($r1.Find.Execute()) {
    $r1.Text = "happened to the mouse when you fan it; use the"
}

# 5. "terminology popular with TV weatherpersons throughout the nation."
#    -> "type terminology popular with TV weatherpersons throughout the nation."
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Text = "terminology popular with TV weatherpersons"
$r2.Find.Forward = $true
if ($r2.Find.Execute()) {
    $r2.Text = "type terminology popular with TV weatherpersons"
}
